$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp label (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Octubre de 2020 a las 16:50"

# --- Afganistan and Birmania swapped rank: row 78 is now Birmania, row 79 is now Afganistan ---
$ws.Range("A78").Value = "Birmania"
$ws.Range("A79").Value = "Afganistan"

# --- Updated case/recovery/death counters for the affected countries ---
$ws.Cells.Item(4,2).Value = 8588169
$ws.Cells.Item(4,3).Value = 3350
$ws.Cells.Item(4,4).Value = 5603138
$ws.Cells.Item(4,5).Value = 2757591
$ws.Cells.Item(4,7).Value = 31
$ws.Cells.Item(4,8).Value = 227440

$ws.Cells.Item(28,2).Value = 307965
$ws.Cells.Item(28,3).Value = 630
$ws.Cells.Item(28,4).Value = 286438
$ws.Cells.Item(28,5).Value = 19208
$ws.Cells.Item(28,7).Value = 28
$ws.Cells.Item(28,8).Value = 2319

$ws.Cells.Item(49,2).Value = 103172
$ws.Cells.Item(49,3).Value = 757
$ws.Cells.Item(49,4).Value = 92665
$ws.Cells.Item(49,5).Value = 6927
$ws.Cells.Item(49,7).Value = 13
$ws.Cells.Item(49,8).Value = 3580

$ws.Cells.Item(52,2).Value = 94524
$ws.Cells.Item(52,3).Value = 591
$ws.Cells.Item(52,4).Value = 87666
$ws.Cells.Item(52,5).Value = 5173
$ws.Cells.Item(52,7).Value = 6
$ws.Cells.Item(52,8).Value = 1685

$ws.Cells.Item(61,2).Value = 69568
$ws.Cells.Item(61,3).Value = 777
$ws.Cells.Item(61,4).Value = 50422
$ws.Cells.Item(61,5).Value = 17505
$ws.Cells.Item(61,7).Value = 11
$ws.Cells.Item(61,8).Value = 1641

$ws.Cells.Item(65,4).Value = 57829
$ws.Cells.Item(65,5).Value = 84

$ws.Cells.Item(78,2).Value = 41008
$ws.Cells.Item(78,3).Value = 1312
$ws.Cells.Item(78,4).Value = 21144
$ws.Cells.Item(78,5).Value = 18859
$ws.Cells.Item(78,7).Value = 33
$ws.Cells.Item(78,8).Value = 1005

$ws.Cells.Item(79,2).Value = 40626
$ws.Cells.Item(79,3).Value = 116
$ws.Cells.Item(79,4).Value = 33831
$ws.Cells.Item(79,5).Value = 5290
$ws.Cells.Item(79,7).Value = 4
$ws.Cells.Item(79,8).Value = 1505

$ws.Cells.Item(95,2).Value = 18250
$ws.Cells.Item(95,3).Value = 302
$ws.Cells.Item(95,4).Value = 10395
$ws.Cells.Item(95,5).Value = 7390
$ws.Cells.Item(95,7).Value = 3
$ws.Cells.Item(95,8).Value = 465

$ws.Cells.Item(97,2).Value = 17125
$ws.Cells.Item(97,3).Value = 161
$ws.Cells.Item(97,5).Value = 4983

$ws.Cells.Item(110,2).Value = 11041
$ws.Cells.Item(110,3).Value = 108
$ws.Cells.Item(110,4).Value = 7210
$ws.Cells.Item(110,5).Value = 3733

$ws.Cells.Item(111,2).Value = 10653
$ws.Cells.Item(111,3).Value = 40
$ws.Cells.Item(111,4).Value = 9724
$ws.Cells.Item(111,5).Value = 848
$ws.Cells.Item(111,7).Value = 1
$ws.Cells.Item(111,8).Value = 81

$ws.Cells.Item(116,2).Value = 8600
$ws.Cells.Item(116,3).Value = 155
$ws.Cells.Item(116,4).Value = 4095
$ws.Cells.Item(116,5).Value = 4326
$ws.Cells.Item(116,7).Value = 5
$ws.Cells.Item(116,8).Value = 179

$ws.Cells.Item(177,2).Value = 551
$ws.Cells.Item(177,3).Value = 1
$ws.Cells.Item(177,5).Value = 53

$ws.Cells.Item(183,2).Value = 457
$ws.Cells.Item(183,3).Value = 5
$ws.Cells.Item(183,4).Value = 391
$ws.Cells.Item(183,5).Value = 66

Write-Host "Edit complete"
